$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.520.25"
$ws.Range("E2").Value = "  +3.52%  "
$ws.Range("D3").Value = "2.422.77"
$ws.Range("E3").Value = "  +1.31%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'563.57"
$ws.Range("E5").Value = "  +2.89%  "
$ws.Range("D6").Value = "'166.32"
$ws.Range("E6").Value = "  +6.24%  "
$ws.Range("D8").Value = "'0.513"
$ws.Range("E8").Value = "  +2.32%  "
$ws.Range("E9").Value = "  +9.89%  "
$ws.Range("D10").Value = "2.421.52"
$ws.Range("E10").Value = "  +1.29%  "
$ws.Range("D12").Value = "'0.334"
$ws.Range("E12").Value = "  +3.40%  "
$ws.Range("D13").Value = "'4.66"
$ws.Range("E13").Value = "  -1.50%  "
$ws.Range("D14").Value = "69.418.07"
$ws.Range("E14").Value = "  +3.44%  "
$ws.Range("E15").Value = "  +6.37%  "
$ws.Range("D16").Value = "2.871.53"
$ws.Range("E16").Value = "  -1.18%  "
$ws.Range("D17").Value = "'23.97"
$ws.Range("E17").Value = "  +5.72%  "
$ws.Range("D18").Value = "2.437.99"
$ws.Range("E18").Value = "  +0.97%  "
$ws.Range("E19").Value = "  +5.46%  "
$ws.Range("D20").Value = "'342.84"
$ws.Range("E20").Value = "  +5.05%  "
$ws.Range("D21").Value = "'7.14"
$ws.Range("E21").Value = "  +5.90%  "
$ws.Range("E22").Value = "  +3.91%  "
$ws.Range("E23").Value = "  +7.67%  "
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("E25").Value = "  +1.18%  "
$ws.Range("D26").Value = "'3.82"
$ws.Range("E26").Value = "  +7.41%  "
$ws.Range("D27").Value = "'8.49"
$ws.Range("E27").Value = "  +7.13%  "
$ws.Range("D28").Value = "2.551.60"
$ws.Range("E28").Value = "  +1.28%  "
$ws.Range("D29").Value = "'0.996"
$ws.Range("E29").Value = "  -0.42%  "
$ws.Range("E30").Value = "  +7.77%  "
$ws.Range("E31").Value = "  +6.79%  "
$ws.Range("E32").Value = "  +12.35%  "
$ws.Range("D33").Value = "'454.05"
$ws.Range("E33").Value = "  +10.28%  "
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("E35").Value = "  +2.41%  "
$ws.Range("D36").Value = "'158.15"
$ws.Range("E36").Value = "  -0.65%  "
$ws.Range("E37").Value = "  +1.01%  "
$ws.Range("E38").Value = "  +6.73%  "
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("D40").Value = "'18.20"
$ws.Range("E40").Value = "  +3.41%  "
$ws.Range("E41").Value = "  +4.58%  "
$ws.Range("E42").Value = "  +5.93%  "
$ws.Range("D43").Value = "'4.40"
$ws.Range("E43").Value = "  +5.64%  "
$ws.Range("D44").Value = "'37.81"
$ws.Range("E44").Value = "  +1.45%  "
$ws.Range("D45").Value = "'1.08"
$ws.Range("E45").Value = "  +3.17%  "
$ws.Range("D46").Value = "'2.08"
$ws.Range("E46").Value = "  +6.86%  "
$ws.Range("D47").Value = "'135.04"
$ws.Range("E47").Value = "  +5.72%  "
$ws.Range("E48").Value = "  +4.09%  "
$ws.Range("E49").Value = "  +3.01%  "
$ws.Range("E50").Value = "  +4.35%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "'0.561"
$ws.Range("E51").Value = "  +2.43%  "
